$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6955149999999999
$ws.Range("H2").Value = 2.086545
$ws.Range("I2").Value = 0.04699474681458952
$ws.Range("J2").Value = 0.04699474681458952
$ws.Range("M2").Value = 32.51511900000001
$ws.Range("N2").Value = 97.54535700000001
$ws.Range("O2").Value = 0.218203973858649
$ws.Range("P2").Value = 0.2182039738586489
$ws.Range("Q2").Value = 22.614752991285
$ws.Range("R2").Value = 203.532776921565
$ws.Range("S2").Value = 0.01025444050542452
$ws.Range("T2").Value = 0.01025444050542452
$ws.Range("G3").Value = 0.6955149999999999
$ws.Range("H3").Value = 2.086545
$ws.Range("I3").Value = 0.04699474681458952
$ws.Range("J3").Value = 0.04699474681458952
$ws.Range("O3").Value = 0.6017421411306194
$ws.Range("P3").Value = 0.6017421411306194
$ws.Range("Q3").Value = 62.36481236098499
$ws.Range("R3").Value = 561.2833112488649
$ws.Range("S3").Value = 0.02827871957010245
$ws.Range("T3").Value = 0.02827871957010245
$ws.Range("G4").Value = 0.6955149999999999
$ws.Range("H4").Value = 2.086545
$ws.Range("I4").Value = 0.04699474681458952
$ws.Range("J4").Value = 0.04699474681458952
$ws.Range("M4").Value = 26.830279
$ws.Range("N4").Value = 80.490837
$ws.Range("O4").Value = 0.1800538850107317
$ws.Range("P4").Value = 0.1800538850107317
$ws.Range("Q4").Value = 18.660861498685
$ws.Range("R4").Value = 167.947753488165
$ws.Range("S4").Value = 0.008461586739062551
$ws.Range("T4").Value = 0.008461586739062549
$ws.Range("I5").Value = 0.8773676626402658
$ws.Range("J5").Value = 0.8773676626402658
$ws.Range("M5").Value = 32.51511900000001
$ws.Range("N5").Value = 97.54535700000001
$ws.Range("O5").Value = 0.218203973858649
$ws.Range("P5").Value = 0.2182039738586489
$ws.Range("Q5").Value = 422.2057637938141
$ws.Range("R5").Value = 3799.851874144326
$ws.Range("S5").Value = 0.1914451105231805
$ws.Range("T5").Value = 0.1914451105231805
$ws.Range("I6").Value = 0.8773676626402658
$ws.Range("J6").Value = 0.8773676626402658
$ws.Range("O6").Value = 0.6017421411306194
$ws.Range("P6").Value = 0.6017421411306194
$ws.Range("S6").Value = 0.5279490958759205
$ws.Range("T6").Value = 0.5279490958759205
$ws.Range("I7").Value = 0.8773676626402658
$ws.Range("J7").Value = 0.8773676626402658
$ws.Range("M7").Value = 26.830279
$ws.Range("N7").Value = 80.490837
$ws.Range("O7").Value = 0.1800538850107317
$ws.Range("P7").Value = 0.1800538850107317
$ws.Range("Q7").Value = 348.388650768774
$ws.Range("R7").Value = 3135.497856918966
$ws.Range("S7").Value = 0.1579734562411649
$ws.Range("T7").Value = 0.1579734562411648
$ws.Range("G8").Value = 1.119424666666667
$ws.Range("H8").Value = 3.358274
$ws.Range("I8").Value = 0.07563759054514464
$ws.Range("J8").Value = 0.07563759054514466
$ws.Range("M8").Value = 32.51511900000001
$ws.Range("N8").Value = 97.54535700000001
$ws.Range("O8").Value = 0.218203973858649
$ws.Range("P8").Value = 0.2182039738586489
$ws.Range("Q8").Value = 36.39822624820201
$ws.Range("R8").Value = 327.584036233818
$ws.Range("S8").Value = 0.01650442283004393
$ws.Range("T8").Value = 0.01650442283004393
$ws.Range("G9").Value = 1.119424666666667
$ws.Range("H9").Value = 3.358274
$ws.Range("I9").Value = 0.07563759054514464
$ws.Range("J9").Value = 0.07563759054514466
$ws.Range("O9").Value = 0.6017421411306194
$ws.Range("P9").Value = 0.6017421411306194
$ws.Range("Q9").Value = 100.375562409042
$ws.Range("R9").Value = 903.3800616813782
$ws.Range("S9").Value = 0.04551432568459643
$ws.Range("T9").Value = 0.04551432568459644
$ws.Range("G10").Value = 1.119424666666667
$ws.Range("H10").Value = 3.358274
$ws.Range("I10").Value = 0.07563759054514464
$ws.Range("J10").Value = 0.07563759054514466
$ws.Range("M10").Value = 26.830279
$ws.Range("N10").Value = 80.490837
$ws.Range("O10").Value = 0.1800538850107317
$ws.Range("P10").Value = 0.1800538850107317
$ws.Range("Q10").Value = 30.03447612614867
$ws.Range("R10").Value = 270.310285135338
$ws.Range("S10").Value = 0.01361884203050428
$ws.Range("T10").Value = 0.01361884203050428
